# BF: Remove (0, 0) as possible position in eyetracking demo
#
# The trials_params sheet lists, for every (target_color, distractor_color)
# ordered pair, one row per possible n_distractors value. Previously
# n_distractors ranged 0..8 (9 rows per color-pair); now it only ranges
# 0..7 (8 rows per color-pair) - i.e. the "8 distractors" row is dropped
# for each pair. This shrinks the sheet from 18 data rows (+1 header) to
# 16 data rows (+1 header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetColor = "black"
$distractorColor = "red"

$row = 2
for ($n = 7; $n -ge 0; $n--) {
    $ws.Cells.Item($row, 1).Value = $targetColor
    $ws.Cells.Item($row, 2).Value = $distractorColor
    $ws.Cells.Item($row, 3).Value = $n
    $row++
}

$targetColor = "red"
$distractorColor = "black"

for ($n = 7; $n -ge 0; $n--) {
    $ws.Cells.Item($row, 1).Value = $targetColor
    $ws.Cells.Item($row, 2).Value = $distractorColor
    $ws.Cells.Item($row, 3).Value = $n
    $row++
}

# Old sheet went up to row 19; the two now-unused trailing rows (18, 19)
# must be cleared out so the used range shrinks to A1:C17.
$ws.Range("A18:C19").Delete()

# Match the author's final selection/cursor position from the diff.
$ws.Range("D9").Select()
